$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D1_USD")
$ws.Range("A1").Value = "test"
